$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New defect rows (B: date, C: title, D: description, E: QC Activity status)
$defectDate = 40837

$rows = @(
    @{ Row = 5;  Title = "Media Option_Permission"; Desc = "Phân quyền Media cho user, dư thừa không cần thiết." },
    @{ Row = 6;  Title = "Media Option_Component"; Desc = "Cấu hình media, dư thừa không cần thiết." },
    @{ Row = 7;  Title = "Category và section"; Desc = "Lỗi giá trị không đồng nhất (giá trị liên kết giữa 2 bảng)" },
    @{ Row = 8;  Title = "Article"; Desc = "Thiếu catid" },
    @{ Row = 9;  Title = "Event trong Article"; Desc = "Thiếu event select change của combox section" },
    @{ Row = 10; Title = "Media Manager"; Desc = "chkThumbnailView và chkDetailView không thể dùng checkbox" },
    @{ Row = 11; Title = "Media Manager"; Desc = "txtFilePath bị dư" },
    @{ Row = 12; Title = "Media Manager"; Desc = "chkSelectFile không cần thiết" },
    @{ Row = 13; Title = "Media Manager"; Desc = "Thiếu sự kiện select chọn trong folder" }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("B$i").Value = $defectDate
    $ws.Range("C$i").Value = $r.Title
    $ws.Range("D$i").Value = $r.Desc
    $ws.Range("E$i").Value = "Error"
}

# Header cell C2 gets a new label (added after the rows so the shared-string
# table ordering matches: row data strings first, header label last)
$ws.Range("C2").Value = "CMS Click and Change"

# Selection moved from G10 to C2
$ws.Range("C2").Select() | Out-Null
